$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records need to be inserted right after the existing row 100
# (i.e. at rows 101-102), which pushes the old rows 101-115 down to 103-117.
$ws.Rows.Item(101).Insert()
$ws.Rows.Item(101).Insert()

# New row 101: Arandano (blue) / Primera / Region de O'Higgins, week of 2021-11-08
$ws.Range("A101").Value = 9
$ws.Range("B101").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C101").Value = "Metropolitana"
$ws.Range("D101").Value = 44508
$ws.Range("E101").Value = 13
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100101
$ws.Range("H101").Value = "Berries"
$ws.Range("I101").Value = 100101001
$ws.Range("J101").Value = "Arándano (blue)"
$ws.Range("K101").Value = "Sin especificar"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 300
$ws.Range("N101").Value = 14000
$ws.Range("O101").Value = 14000
$ws.Range("P101").Value = 14000
$ws.Range("Q101").Value = "$/bandeja 2 kilos"
$ws.Range("R101").Value = "Región de O'Higgins"
$ws.Range("S101").Value = 7000
$ws.Range("T101").Value = 2

# New row 102: Arandano (blue) / Segunda / Region de O'Higgins, week of 2021-11-08
$ws.Range("A102").Value = 9
$ws.Range("B102").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C102").Value = "Metropolitana"
$ws.Range("D102").Value = 44508
$ws.Range("E102").Value = 13
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100101
$ws.Range("H102").Value = "Berries"
$ws.Range("I102").Value = 100101001
$ws.Range("J102").Value = "Arándano (blue)"
$ws.Range("K102").Value = "Sin especificar"
$ws.Range("L102").Value = "Segunda"
$ws.Range("M102").Value = 280
$ws.Range("N102").Value = 12000
$ws.Range("O102").Value = 12000
$ws.Range("P102").Value = 12000
$ws.Range("Q102").Value = "$/bandeja 2 kilos"
$ws.Range("R102").Value = "Región de O'Higgins"
$ws.Range("S102").Value = 6000
$ws.Range("T102").Value = 2
